# Update for 2017-01-31 release: extend Table 6.7.B through November 2016.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the report title (row 1) from "October 2016" to "November 2016".
$ws.Range("A1").Value = "Table 6.7.B. Capacity Factors for Utility Scale Generators Not Primarily Using Fossil Fuels, January 2013-November 2016"

# 2) Insert a new data row for "November" right before the footnote row
#    (old row 44), which pushes the footnote row down to row 45.
$ws.Rows.Item(44).Insert()

# Copy formatting from the row above (October, row 43) into the new row
# so the new row's styles match the rest of the monthly data rows.
$ws.Range("A43:I43").Copy()
$ws.Range("A44:I44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the November values.
$ws.Range("A44").Value = "November"
$ws.Range("B44").Value = 0.911
$ws.Range("C44").Value = 0.333
$ws.Range("D44").Value = 0.355
$ws.Range("E44").Value = 0.213
$ws.Range("F44").Value = 0.144
$ws.Range("G44").Value = 0.715
$ws.Range("H44").Value = 0.416
$ws.Range("I44").Value = 0.782
